$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Jasmine Day"
$ws.Range("B2").Value = "3************2"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "28/09/1933"
$ws.Range("C2").ClearFormats()
$ws.Range("D2").Value = "95604 Cassandra Road Apt. 633`nPort Victoriafort, CO 71069"
$ws.Range("F2").Value = "Perform address."
$ws.Range("H2").Value = "North Jacob"
$ws.Range("I2").Value = "No up citizen when."
$ws.Range("J2").Value = 47
$ws.Range("K2").Value = "State including."
$ws.Range("L2").Value = "Determine her."
$ws.Range("M2").Value = 3948
$ws.Range("N2").Value = 9633
$ws.Range("O2").Value = "shake"

# Row 3
$ws.Range("A3").Value = "Chris Villegas"
$ws.Range("B3").Value = "8************5"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "14/01/1964"
$ws.Range("C3").ClearFormats()
$ws.Range("D3").Value = "356 Townsend Islands`nLake Matthew, NJ 15801"
$ws.Range("F3").Value = "Theory move."
$ws.Range("H3").Value = "East Henry"
$ws.Range("I3").Value = "Enter sometimes PM."
$ws.Range("J3").Value = 90
$ws.Range("K3").Value = "Born political use."
$ws.Range("L3").Value = "Mention before list."
$ws.Range("M3").Value = 2016
$ws.Range("N3").Value = 3598
$ws.Range("O3").Value = "offer"

# Row 4
$ws.Range("A4").Value = "Donna Walsh"
$ws.Range("B4").Value = "8************6"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "16/09/1976"
$ws.Range("C4").ClearFormats()
$ws.Range("D4").Value = "5040 Bryant Motorway`nWest Janet, IA 23729"
$ws.Range("F4").Value = "Nor rich debate."
$ws.Range("H4").Value = "North Codychester"
$ws.Range("I4").Value = "Body treatment."
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = "Subject human."
$ws.Range("L4").Value = "Yeah when product."
$ws.Range("M4").Value = 8352
$ws.Range("N4").Value = 4204
$ws.Range("O4").Value = "animal"

# Row 5
$ws.Range("A5").Value = "Ian Vasquez DDS"
$ws.Range("B5").Value = "4************8"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "01/08/1933"
$ws.Range("C5").ClearFormats()
$ws.Range("D5").Value = "233 Zachary Course Apt. 717`nPerryborough, WA 82106"
$ws.Range("F5").Value = "Over position."
$ws.Range("H5").Value = "Matthewhaven"
$ws.Range("I5").Value = "Of room measure."
$ws.Range("J5").Value = 19
$ws.Range("K5").Value = "Congress move begin."
$ws.Range("L5").Value = "No method start."
$ws.Range("M5").Value = 2642
$ws.Range("N5").Value = 6425
$ws.Range("O5").Value = "work"

# Row 6
$ws.Range("A6").Value = "Alex Russell"
$ws.Range("B6").Value = "6************6"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "18/01/1927"
$ws.Range("C6").ClearFormats()
$ws.Range("D6").Value = "144 Brett Park Suite 101`nLake Christopher, MD 55636"
$ws.Range("F6").Value = "Officer yourself."
$ws.Range("H6").Value = "Gonzalesmouth"
$ws.Range("I6").Value = "Beyond language."
$ws.Range("J6").Value = 21
$ws.Range("K6").Value = "By rich training."
$ws.Range("L6").Value = "Cost should second."
$ws.Range("M6").Value = 4940
$ws.Range("N6").Value = 9056
$ws.Range("O6").Value = "will"

# Row 7
$ws.Range("A7").Value = "Jason Elliott"
$ws.Range("B7").Value = "7************8"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "02/06/1947"
$ws.Range("C7").ClearFormats()
$ws.Range("D7").Value = "USCGC Carlson`nFPO AA 19018"
$ws.Range("F7").Value = "News social list."
$ws.Range("H7").Value = "Lake Matthewberg"
$ws.Range("I7").Value = "Only maybe history."
$ws.Range("J7").Value = 38
$ws.Range("K7").Value = "Wrong our article."
$ws.Range("L7").Value = "Local strategy."
$ws.Range("M7").Value = 3499
$ws.Range("N7").Value = 7904
$ws.Range("O7").Value = "live"

# Row 8
$ws.Range("A8").Value = "John Collins"
$ws.Range("B8").Value = "8************0"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "31/08/1933"
$ws.Range("C8").ClearFormats()
$ws.Range("D8").Value = "2189 Martin Street`nNorth Amanda, GU 59555"
$ws.Range("F8").Value = "Various oil what."
$ws.Range("H8").Value = "Carpenterfurt"
$ws.Range("I8").Value = "Future able street."
$ws.Range("J8").Value = 57
$ws.Range("K8").Value = "Can difference term."
$ws.Range("L8").Value = "Father decide key."
$ws.Range("M8").Value = 8112
$ws.Range("N8").Value = 7459
$ws.Range("O8").Value = "within"

# Row 9
$ws.Range("A9").Value = "Erica Keller"
$ws.Range("B9").Value = "9************6"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "20/08/1985"
$ws.Range("C9").ClearFormats()
$ws.Range("D9").Value = "853 Allen Plains Apt. 968`nHowardborough, IA 63954"
$ws.Range("F9").Value = "Daughter cause."
$ws.Range("H9").Value = "East Amy"
$ws.Range("I9").Value = "Room toward before."
$ws.Range("J9").Value = 98
$ws.Range("K9").Value = "Bag between leader."
$ws.Range("L9").Value = "Business economic."
$ws.Range("M9").Value = 2720
$ws.Range("N9").Value = 5999
$ws.Range("O9").Value = "behind"

# Row 10
$ws.Range("A10").Value = "Jennifer Butler"
$ws.Range("B10").Value = "2************0"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "17/12/1972"
$ws.Range("C10").ClearFormats()
$ws.Range("D10").Value = "1259 Baker Parks Suite 277`nBurkeland, CO 62870"
$ws.Range("F10").Value = "Success green."
$ws.Range("H10").Value = "Campbellfort"
$ws.Range("I10").Value = "Ahead them fall."
$ws.Range("J10").Value = 81
$ws.Range("K10").Value = "Factor him cause."
$ws.Range("L10").Value = "Me camera shake."
$ws.Range("M10").Value = 4319
$ws.Range("N10").Value = 3134
$ws.Range("O10").Value = "building"

# Row 11
$ws.Range("A11").Value = "Linda Tran"
$ws.Range("B11").Value = "5************8"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "02/01/1946"
$ws.Range("C11").ClearFormats()
$ws.Range("D11").Value = "7114 Sanchez Canyon Suite 597`nHoustontown, IN 40676"
$ws.Range("F11").Value = "Training attorney."
$ws.Range("H11").Value = "Kennethmouth"
$ws.Range("I11").Value = "Officer me speak."
$ws.Range("J11").Value = 62
$ws.Range("K11").Value = "Much name protect."
$ws.Range("L11").Value = "Show moment."
$ws.Range("M11").Value = 3825
$ws.Range("N11").Value = 6888
$ws.Range("O11").Value = "tree"

# Row 12
$ws.Range("A12").Value = "Jessica Brown"
$ws.Range("B12").Value = "4************1"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "19/02/1989"
$ws.Range("C12").ClearFormats()
$ws.Range("D12").Value = "852 Darrell Estate`nAndrewborough, MN 37579"
$ws.Range("F12").Value = "Pull husband."
$ws.Range("H12").Value = "East Megan"
$ws.Range("I12").Value = "Space tell suffer."
$ws.Range("J12").Value = 61
$ws.Range("K12").Value = "Watch position."
$ws.Range("L12").Value = "Human theory coach."
$ws.Range("M12").Value = 6394
$ws.Range("N12").Value = 6384
$ws.Range("O12").Value = "spend"

# Row 13
$ws.Range("A13").Value = "Tyler Wilkerson"
$ws.Range("B13").Value = "0************7"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "20/01/1936"
$ws.Range("C13").ClearFormats()
$ws.Range("D13").Value = "9455 Jeffrey Court Suite 350`nNorth Susanfurt, CT 30635"
$ws.Range("F13").Value = "Relationship hit."
$ws.Range("H13").Value = "South Carlland"
$ws.Range("I13").Value = "Billion real."
$ws.Range("J13").Value = 81
$ws.Range("K13").Value = "Direction prove."
$ws.Range("L13").Value = "One dark democratic."
$ws.Range("M13").Value = 9039
$ws.Range("N13").Value = 2401
$ws.Range("O13").Value = "majority"

# Row 14
$ws.Range("A14").Value = "Taylor Rowe"
$ws.Range("B14").Value = "5************3"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "27/02/1974"
$ws.Range("C14").ClearFormats()
$ws.Range("D14").Value = "8452 Bennett River Suite 609`nWest Heather, SC 96601"
$ws.Range("F14").Value = "Record social every."
$ws.Range("H14").Value = "Davidburgh"
$ws.Range("I14").Value = "Thousand oil sense."
$ws.Range("J14").Value = 30
$ws.Range("K14").Value = "Attorney recent."
$ws.Range("L14").Value = "Before issue event."
$ws.Range("M14").Value = 6887
$ws.Range("N14").Value = 9577
$ws.Range("O14").Value = "blue"

# Row 15
$ws.Range("A15").Value = "Lindsey Stuart"
$ws.Range("B15").Value = "8************2"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "27/04/1949"
$ws.Range("C15").ClearFormats()
$ws.Range("D15").Value = "55064 Munoz Meadows`nWest Autumn, TX 49288"
$ws.Range("F15").Value = "Away voice."
$ws.Range("H15").Value = "South Teresastad"
$ws.Range("I15").Value = "Guess unit heart."
$ws.Range("J15").Value = 63
$ws.Range("K15").Value = "Same great still."
$ws.Range("L15").Value = "To under movie ask."
$ws.Range("M15").Value = 5929
$ws.Range("N15").Value = 6438
$ws.Range("O15").Value = "group"

# Row 16
$ws.Range("A16").Value = "Michael Butler"
$ws.Range("B16").Value = "5************1"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "12/12/1977"
$ws.Range("C16").ClearFormats()
$ws.Range("D16").Value = "26513 Campbell Stravenue Apt. 183`nWest Scottfort, UT 28555"
$ws.Range("F16").Value = "Outside best."
$ws.Range("H16").Value = "Tiffanyview"
$ws.Range("I16").Value = "Week part thank."
$ws.Range("J16").Value = 28
$ws.Range("K16").Value = "Rich scientist."
$ws.Range("L16").Value = "Manager almost."
$ws.Range("M16").Value = 5419
$ws.Range("N16").Value = 6910
$ws.Range("O16").Value = "role"

# Row 17
$ws.Range("A17").Value = "Vanessa Campos"
$ws.Range("B17").Value = "3************1"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "24/07/1982"
$ws.Range("C17").ClearFormats()
$ws.Range("D17").Value = "86212 Johnson Ways Apt. 892`nSmithfort, ID 90278"
$ws.Range("F17").Value = "Girl require writer."
$ws.Range("H17").Value = "East James"
$ws.Range("I17").Value = "None speech."
$ws.Range("J17").Value = 53
$ws.Range("K17").Value = "Position best that."
$ws.Range("L17").Value = "Ground teach short."
$ws.Range("M17").Value = 4849
$ws.Range("N17").Value = 8887
$ws.Range("O17").Value = "art"

# Row 18
$ws.Range("A18").Value = "Timothy Moreno"
$ws.Range("B18").Value = "7************9"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "24/05/2006"
$ws.Range("C18").ClearFormats()
$ws.Range("D18").Value = "7518 Barbara Manors`nTimothymouth, OH 79682"
$ws.Range("F18").Value = "Position trip."
$ws.Range("H18").Value = "West Stevenstad"
$ws.Range("I18").Value = "Activity thank."
$ws.Range("J18").Value = 48
$ws.Range("K18").Value = "Money kind reflect."
$ws.Range("L18").Value = "Others work discuss."
$ws.Range("M18").Value = 8809
$ws.Range("N18").Value = 9891
$ws.Range("O18").Value = "offer"

# Row 19
$ws.Range("A19").Value = "Ashley Roach"
$ws.Range("B19").Value = "9************5"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "17/04/1957"
$ws.Range("C19").ClearFormats()
$ws.Range("D19").Value = "0106 Henry Stravenue`nEast Angela, NC 74149"
$ws.Range("F19").Value = "Wear growth system."
$ws.Range("H19").Value = "Karenhaven"
$ws.Range("I19").Value = "Rate kitchen fish."
$ws.Range("J19").Value = 45
$ws.Range("K19").Value = "Rate meeting chair."
$ws.Range("L19").Value = "Quickly although."
$ws.Range("M19").Value = 2166
$ws.Range("N19").Value = 4469
$ws.Range("O19").Value = "pretty"

# Row 20
$ws.Range("A20").Value = "Jenna Jones"
$ws.Range("B20").Value = "4************6"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "23/07/1935"
$ws.Range("C20").ClearFormats()
$ws.Range("D20").Value = "2895 Pacheco Ranch Apt. 824`nSimsborough, PA 31956"
$ws.Range("F20").Value = "Remember give."
$ws.Range("H20").Value = "North Donna"
$ws.Range("I20").Value = "Make into card firm."
$ws.Range("J20").Value = 37
$ws.Range("K20").Value = "Part wide current."
$ws.Range("L20").Value = "Pretty key."
$ws.Range("M20").Value = 9836
$ws.Range("N20").Value = 6315
$ws.Range("O20").Value = "remember"

# Row 21
$ws.Range("A21").Value = "Ashley Hansen"
$ws.Range("B21").Value = "4************4"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "23/06/1927"
$ws.Range("C21").ClearFormats()
$ws.Range("D21").Value = "USS Garcia`nFPO AA 51097"
$ws.Range("F21").Value = "View join walk cell."
$ws.Range("H21").Value = "Byrdton"
$ws.Range("I21").Value = "Second lot if hold."
$ws.Range("J21").Value = 82
$ws.Range("K21").Value = "Science management."
$ws.Range("L21").Value = "Interest perform."
$ws.Range("M21").Value = 9461
$ws.Range("N21").Value = 4164
$ws.Range("O21").Value = "daughter"

# Row 22
$ws.Range("A22").Value = "Emma Marshall"
$ws.Range("B22").Value = "4************3"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "16/03/1944"
$ws.Range("C22").ClearFormats()
$ws.Range("D22").Value = "2579 Hall Lights`nJosephburgh, IN 92388"
$ws.Range("F22").Value = "Daughter must upon."
$ws.Range("H22").Value = "Andrewhaven"
$ws.Range("I22").Value = "Mean avoid itself."
$ws.Range("J22").Value = 48
$ws.Range("K22").Value = "Pay among human."
$ws.Range("L22").Value = "Woman budget sport."
$ws.Range("M22").Value = 6664
$ws.Range("N22").Value = 4963
$ws.Range("O22").Value = "firm"

# Row 23
$ws.Range("A23").Value = "Harry Wilson"
$ws.Range("B23").Value = "0************0"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "09/12/2005"
$ws.Range("C23").ClearFormats()
$ws.Range("D23").Value = "524 Jonathan Causeway Apt. 771`nAnnehaven, MH 35480"
$ws.Range("F23").Value = "Cut who raise."
$ws.Range("H23").Value = "North Brittany"
$ws.Range("I23").Value = "Say decision find."
$ws.Range("J23").Value = 70
$ws.Range("K23").Value = "For surface stop."
$ws.Range("L23").Value = "Evening value him."
$ws.Range("M23").Value = 3063
$ws.Range("N23").Value = 8773
$ws.Range("O23").Value = "yeah"

# Row 24
$ws.Range("A24").Value = "Patricia Kramer"
$ws.Range("B24").Value = "7************0"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "24/10/1979"
$ws.Range("C24").ClearFormats()
$ws.Range("D24").Value = "78667 Jeffrey Grove`nLake Marissachester, RI 23908"
$ws.Range("F24").Value = "Mean total state."
$ws.Range("H24").Value = "North Coleshire"
$ws.Range("I24").Value = "Through try skill."
$ws.Range("J24").Value = 29
$ws.Range("K24").Value = "Subject record."
$ws.Range("L24").Value = "Page wish meet my."
$ws.Range("M24").Value = 6323
$ws.Range("N24").Value = 6877
$ws.Range("O24").Value = "shoulder"

# Row 25
$ws.Range("A25").Value = "Kaitlin Brennan"
$ws.Range("B25").Value = "0************2"
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = "13/08/1975"
$ws.Range("C25").ClearFormats()
$ws.Range("D25").Value = "973 Gardner Hollow`nBallchester, FM 35574"
$ws.Range("F25").Value = "Agent nice voice."
$ws.Range("H25").Value = "East Frank"
$ws.Range("I25").Value = "Forget imagine."
$ws.Range("J25").Value = 70
$ws.Range("K25").Value = "Laugh industry."
$ws.Range("L25").Value = "Happy standard."
$ws.Range("M25").Value = 8758
$ws.Range("N25").Value = 8357
$ws.Range("O25").Value = "process"

# Row 26
$ws.Range("A26").Value = "Molly Mcknight"
$ws.Range("B26").Value = "7************4"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "31/08/1966"
$ws.Range("C26").ClearFormats()
$ws.Range("D26").Value = "27667 Gates Cove Suite 107`nNicolefort, TX 33737"
$ws.Range("F26").Value = "Production sense."
$ws.Range("H26").Value = "New Keith"
$ws.Range("I26").Value = "Recognize challenge."
$ws.Range("J26").Value = 25
$ws.Range("K26").Value = "Outside reflect."
$ws.Range("L26").Value = "Begin one anything."
$ws.Range("M26").Value = 8511
$ws.Range("N26").Value = 7165
$ws.Range("O26").Value = "community"

# Row 27
$ws.Range("A27").Value = "Jennifer Castro"
$ws.Range("B27").Value = "4************0"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "08/02/1974"
$ws.Range("C27").ClearFormats()
$ws.Range("D27").Value = "78506 Nicole Cliff`nPort Kevinborough, TN 06706"
$ws.Range("F27").Value = "Save bank TV."
$ws.Range("H27").Value = "Port Gregory"
$ws.Range("I27").Value = "Social method."
$ws.Range("J27").Value = 43
$ws.Range("K27").Value = "Rich or firm."
$ws.Range("L27").Value = "Anything body."
$ws.Range("M27").Value = 4941
$ws.Range("N27").Value = 5423
$ws.Range("O27").Value = "though"

# Row 28
$ws.Range("A28").Value = "Cheryl Lee"
$ws.Range("B28").Value = "4************0"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "19/05/1972"
$ws.Range("C28").ClearFormats()
$ws.Range("D28").Value = "7825 West Forest`nEast Jason, KS 35208"
$ws.Range("F28").Value = "Simple fine often."
$ws.Range("H28").Value = "Lindsayview"
$ws.Range("I28").Value = "Share listen."
$ws.Range("J28").Value = 46
$ws.Range("K28").Value = "International need."
$ws.Range("L28").Value = "Lead forget six."
$ws.Range("M28").Value = 3895
$ws.Range("N28").Value = 8207
$ws.Range("O28").Value = "ability"

# Row 29
$ws.Range("A29").Value = "Ryan Brown"
$ws.Range("B29").Value = "5************5"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "17/11/1986"
$ws.Range("C29").ClearFormats()
$ws.Range("D29").Value = "115 Cardenas Mountain Suite 333`nChristinamouth, LA 83337"
$ws.Range("F29").Value = "Play full team."
$ws.Range("H29").Value = "North Coryview"
$ws.Range("I29").Value = "Fill culture read."
$ws.Range("J29").Value = 19
$ws.Range("K29").Value = "Hope into social."
$ws.Range("L29").Value = "Remember option."
$ws.Range("M29").Value = 9140
$ws.Range("N29").Value = 8056
$ws.Range("O29").Value = "accept"

# Row 30
$ws.Range("A30").Value = "Jennifer Heath"
$ws.Range("B30").Value = "2************0"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "09/10/1935"
$ws.Range("C30").ClearFormats()
$ws.Range("D30").Value = "5778 Burgess Stravenue`nPort Patriciaside, PA 48083"
$ws.Range("F30").Value = "Anyone mother."
$ws.Range("H30").Value = "East James"
$ws.Range("I30").Value = "Position challenge."
$ws.Range("J30").Value = 21
$ws.Range("K30").Value = "Chair statement no."
$ws.Range("L30").Value = "Stand federal final."
$ws.Range("M30").Value = 9582
$ws.Range("N30").Value = 7280
$ws.Range("O30").Value = "over"

# Row 31
$ws.Range("A31").Value = "Holly Torres"
$ws.Range("B31").Value = "7************3"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "14/07/1974"
$ws.Range("C31").ClearFormats()
$ws.Range("D31").Value = "PSC 7807, Box 0960`nAPO AA 28266"
$ws.Range("F31").Value = "Large recently let."
$ws.Range("H31").Value = "Kaylaborough"
$ws.Range("I31").Value = "Condition billion."
$ws.Range("J31").Value = 88
$ws.Range("K31").Value = "Think plant space."
$ws.Range("L31").Value = "Black TV everything."
$ws.Range("M31").Value = 7819
$ws.Range("N31").Value = 3449
$ws.Range("O31").Value = "inside"
